$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A5 with the refreshed timestamp value
$ws.Range("A5").Value = 45873.50022493055

# Add new row 6 with the latest sensor reading
$ws.Range("A6").Value = 45873.5419174335
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = 15
$ws.Range("D6").Value = 19.89
$ws.Range("E6").Value = 75.47
$ws.Range("F6").Value = 620.51
$ws.Range("G6").Value = 12.24
$ws.Range("H6").Value = "ESE"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "13:00:21"

# Match the number format of A5 (date formatted) for the new A6 cell
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
